# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The error-table rows were shifted by one quarter (an off-by-one in the
# naive component forecaster). Row 2 (Q0) now holds newly computed error
# stats, and what used to be in rows 2-10 now belongs one row down
# (rows 3-11); the N counts (column G) increase by one accordingly. The
# former last row (row 11 / Q9) of old data is superseded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing B2:G10 block (old rows 2-10) before overwriting
# anything, then shift it down into B3:G11.
$oldBlock = $ws.Range("B2:G10").Value()

$ws.Range("B3:G11").Value = $oldBlock

# Write the newly computed values for row 2 (Q0).
$ws.Range("B2").Value = 0.03864900479865655
$ws.Range("C2").Value = 0.5850438662025559
$ws.Range("D2").Value = 0.7245914535781601
$ws.Range("E2").Value = 0.8512293777696821
$ws.Range("F2").Value = 0.8736531364880852
$ws.Range("G2").Value = 19
